$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.20"
$ws.Range("D3").Value = "'21.75"
$ws.Range("D4").Value = "'5.363"
$ws.Range("D5").Value = "'0.05613"
$ws.Range("D6").Value = "'3.406"
$ws.Range("D7").Value = "'6.384"
$ws.Range("D9").Value = "'0.9517"
$ws.Range("D11").Value = "'0.07628"
$ws.Range("D12").Value = "'0.03181"
$ws.Range("D15").Value = "'3.554"
$ws.Range("D16").Value = "'0.001592"
$ws.Range("D17").Value = "'0.04707"
$ws.Range("D18").Value = "'0.0005767"
$ws.Range("D19").Value = "'0.006273"
$ws.Range("D20").Value = "'0.005081"
$ws.Range("D21").Value = "'0.001032"
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("D23").Value = "'3.749"
$ws.Range("D24").Value = "'2.140"
$ws.Range("D25").Value = "'0.3251"
$ws.Range("D40").Value = "'0.03947"
$ws.Range("D41").Value = "'0.006988"
$ws.Range("D42").Value = "'0.1064"
$ws.Range("D43").Value = "'0.003028"
$ws.Range("D44").Value = "'0.008598"
$ws.Range("D47").Value = "'0.0005497"
$ws.Range("D49").Value = "'0.1642"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D51").Value = "'0.01009"
